# Highlight three task paragraphs in yellow (H2 -> MongoDB related task list
# items now flagged as changed/important):
#   - "Need to return the exercise in a week."
#   - "Persistent the Data on a saas DB."
#   - "Documenting the steps in a README file."
#
# wdYellow = 7 for WdColorIndex / HighlightColorIndex

$d = $word.ActiveDocument

$targets = @(
    "Need to return the exercise in a week.",
    "Persistent the Data on a saas DB.",
    "Documenting the steps in a README file."
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    foreach ($t in $targets) {
        if ($text -eq $t) {
            $p.Range.HighlightColorIndex = 7
        }
    }
}
